# Update vessel activity mapping: insert a new data row (Onshore / Transit)
# and shift the "Towing Group" rows down, re-pointing the remaining values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 5 ("Towing Group" / "Idle at port"),
# shifting it (and the rows below it) down by one row.
$ws.Rows(5).Insert()

# Copy the formatting of the index column from row 4 into the newly
# inserted row 5 so it keeps the same bold/centered style as its neighbours.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)

# --- Row 2: HLV / Idle at port -> Onshore / Idle at port -------------------
$ws.Range("B2").Value = "Onshore"
$ws.Range("C2").Value = "Idle at port"
$ws.Range("D2").Value = 24635.25
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()

# --- Row 3: HLV / Transit -> Onshore / None ---------------------------------
$ws.Range("B3").Value = "Onshore"
$ws.Range("C3").Value = "None"
$ws.Range("D3").Value = 424809
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()

# --- Row 4: Towing Group / Idle at port -> Onshore / Transit ---------------
$ws.Range("B4").Value = "Onshore"
$ws.Range("C4").Value = "Transit"
$ws.Range("D4").Value = 11088
$ws.Range("E4").ClearContents()
$ws.Range("F4").ClearContents()

# --- Row 5 (new): Towing Group / Idle at port -------------------------------
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Towing Group"
$ws.Range("C5").Value = "Idle at port"
$ws.Range("D5").Value = 678265.29262013
$ws.Range("E5").Value = 0.407310331
$ws.Range("F5").Value = 276264.460842917
$ws.Range("G5").Value = "25 GW (SC)"

# --- Row 6: Towing Group / Idle at sea --------------------------------------
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Towing Group"
$ws.Range("C6").Value = "Idle at sea"
$ws.Range("D6").Value = 108108
$ws.Range("E6").Value = 0.203655165
$ws.Range("F6").Value = 22016.75257782
$ws.Range("G6").Value = "25 GW (SC)"

# --- Row 7: Towing Group / Maneuvering --------------------------------------
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Towing Group"
$ws.Range("C7").Value = "Maneuvering"
$ws.Range("D7").Value = 8316
$ws.Range("E7").Value = 2.698858249
$ws.Range("F7").Value = 22443.705198684
$ws.Range("G7").Value = "25 GW (SC)"

# --- Row 8: Towing Group / Transit ------------------------------------------
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Towing Group"
$ws.Range("C8").Value = "Transit"
$ws.Range("D8").Value = 198897.0324545454
$ws.Range("E8").Value = 4.753093345
$ws.Range("F8").Value = 945376.1612999489
$ws.Range("G8").Value = "25 GW (SC)"
